# Insert a new weekly price record for "Ajo" (Chino / Primera) at Terminal
# Hortofrutícola Agro Chillán, pushing the existing historical rows for this
# sub-block (74..162) down by one row (75..163) and adding the new
# observation at row 74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 74..162 down to 75..163, leaving row 74 free for the new entry.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row with the new price observation.
$ws.Range("A74").Value = 7
$ws.Range("B74").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C74").Value = "Ñuble"
$ws.Range("D74").Value = 44539
$ws.Range("E74").Value = 16
$ws.Range("F74").Value = 100112003
$ws.Range("G74").Value = "Ajo"
$ws.Range("H74").Value = "Chino"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 100
$ws.Range("K74").Value = 19000
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = 19500
$ws.Range("N74").Value = "`$/caja 10 kilos"
$ws.Range("O74").Value = "China"
$ws.Range("P74").Value = 1950
$ws.Range("Q74").Value = 10
$ws.Range("R74").Value = "Hortaliza"
